{"js": "// Apply hybrid bold + color (2C3E50) highlighting to quantitative impact\n// metrics (percentages, dollar amounts, large numbers) across the resume's\n// achievements / work-experience bullet paragraphs, matching the target\n// OOXML diff: each metric becomes its own run with <w:b/> and\n// <w:color w:val=\"2C3E50\"/>, while the surrounding text stays in plain runs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map each target paragraph (matched by its exact original text) to the\n// ordered list of metric substrings inside it that must become bold+colored\n// runs. Substrings are unique within their paragraph, so a simple\n// paragraph-scoped search finds exactly one hit each.\nconst targets = [\n  {\n    text:\n      \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    text:\n      \"\\u2022 Utilized advanced sampling methods to decrease survey margin of error from \\u00b14.2% to \\u00b12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\",\n    metrics: [\"\\u00b14.2%\", \"\\u00b12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    text:\n      \"\\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    text:\n      \"\\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n    metrics: [\"$2\"],\n  },\n  {\n    text:\n      \"\\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n    metrics: [\"73.5%\"],\n  },\n  {\n    text: \"\\u2022 $4.7M savings enabled nonprofit access\",\n    metrics: [\"$4.7M\"],\n  },\n  {\n    text: \"\\u2022 178% accuracy improvement in racial classification algorithms\",\n    metrics: [\"178%\"],\n  },\n];\n\nfor (const target of targets) {\n  // Find the index of the paragraph whose full text matches exactly.\n  let para = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === target.text) {\n      para = paragraphs.items[i];\n      break;\n    }\n  }\n  if (!para) {\n    continue; // already edited, or not found \u2014 skip defensively\n  }\n\n  for (const metric of target.metrics) {\n    const found = para.search(metric, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n    if (found.items.length === 0) continue;\n    const range = found.items[0];\n    range.font.bold = true;\n    range.font.color = \"#2C3E50\";\n  }\n  await context.sync();\n}\n", "ps1": "# Apply hybrid bold + color (2C3E50) highlighting to quantitative impact\n# metrics (percentages, dollar amounts, large numbers) across the resume's\n# achievements / work-experience bullet paragraphs, matching the target\n# OOXML diff: each metric becomes its own run with <w:b/> and\n# <w:color w:val=\"2C3E50\"/>, while the surrounding text stays in plain runs.\n\n$d = $word.ActiveDocument\n\n$bullet = [char]0x2022\n$pm = [char]0xb1          # \"\u00b1\"\n$NAVY = 5258796            # RGB(0x2C,0x3E,0x50) packed as BGR long -> 0x503E2C\n\n$targets = @(\n    @{\n        Text    = $bullet + \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Metrics = @(\"23%\", \"64%\")\n    },\n    @{\n        Text    = $bullet + \" Utilized advanced sampling methods to decrease survey margin of error from \" + $pm + \"4.2% to \" + $pm + \"2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\"\n        Metrics = @($pm + \"4.2%\", $pm + \"2.1%\", \"71%\", \"87%\")\n    },\n    @{\n        Text    = $bullet + \" Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\"\n        Metrics = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        Text    = $bullet + \" Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\"\n        Metrics = @(\"`$2\")\n    },\n    @{\n        Text    = $bullet + \" Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\"\n        Metrics = @(\"73.5%\")\n    },\n    @{\n        Text    = $bullet + \" `$4.7M savings enabled nonprofit access\"\n        Metrics = @(\"`$4.7M\")\n    },\n    @{\n        Text    = $bullet + \" 178% accuracy improvement in racial classification algorithms\"\n        Metrics = @(\"178%\")\n    }\n)\n\nforeach ($target in $targets) {\n    $para = $null\n    foreach ($p in $d.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd([char]0x07, [char]0x0d, [char]0x0a)\n        if ($t -eq $target.Text) {\n            $para = $p\n            break\n        }\n    }\n    if ($para -eq $null) {\n        continue\n    }\n\n    $paraStart = $para.Range.Start\n    $paraEnd = $para.Range.End\n\n    foreach ($metric in $target.Metrics) {\n        $r = $d.Range($paraStart, $paraEnd)\n        $r.Find.ClearFormatting()\n        $r.Find.Text = $metric\n        $r.Find.Forward = $true\n        $found = $r.Find.Execute()\n        if ($found) {\n            $r.Font.Bold = 1\n            $r.Font.Color = $NAVY\n        }\n    }\n}\n"}
